$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 29.34999999999999
$ws.Range("C2").Value = 29.78336143493652
$ws.Range("D2").Value = 0.4333614349365291
$ws.Range("E2").Value = 0.1878021332902476
$ws.Range("C3").Value = 29.6746883392334
$ws.Range("D3").Value = 0.3046883392333939
$ws.Range("E3").Value = 0.09283498406480371
$ws.Range("B4").Value = 29.53999999999999
$ws.Range("C4").Value = 29.21548271179199
$ws.Range("D4").Value = -0.3245172882079999
$ws.Range("E4").Value = 0.105311470345874
$ws.Range("C5").Value = 29.43997001647949
$ws.Range("D5").Value = -0.110029983520505
$ws.Range("E5").Value = 0.0121065972735226
$ws.Range("C6").Value = 29.63873863220215
$ws.Range("D6").Value = -0.1112613677978516
$ws.Range("E6").Value = 0.0123790919642488
$ws.Range("C7").Value = 29.90313529968262
$ws.Range("D7").Value = 0.06313529968261378
$ws.Range("E7").Value = 0.003986066066013452
$ws.Range("C8").Value = 29.88266563415527
$ws.Range("D8").Value = 0.07266563415527116
$ws.Range("E8").Value = 0.005280294387187711
$ws.Range("C9").Value = 29.90141868591309
$ws.Range("D9").Value = -0.01858131408691577
$ws.Range("E9").Value = 0.0003452652331966143
$ws.Range("C10").Value = 29.95934104919434
$ws.Range("D10").Value = -0.02065895080566804
$ws.Range("E10").Value = 0.0004267922483910122
$ws.Range("B11").Value = 30.03999999999999
$ws.Range("C11").Value = 30.1280574798584
$ws.Range("D11").Value = 0.0880574798584064
$ws.Range("E11").Value = 0.007754119759013648
$ws.Range("B12").Value = 30.21000000000001
$ws.Range("C12").Value = 30.17116355895996
$ws.Range("D12").Value = -0.03883644104004702
$ws.Range("E12").Value = 0.001508269152657048
$ws.Range("C13").Value = 30.2900447845459
$ws.Range("D13").Value = 0.07004478454589957
$ws.Range("E13").Value = 0.004906271842081492
$ws.Range("C14").Value = 30.34025382995605
$ws.Range("D14").Value = -0.03974617004394077
$ws.Range("E14").Value = 0.001579758033161854
$ws.Range("C15").Value = 30.54611778259277
$ws.Range("D15").Value = 0.1061177825927757
$ws.Range("E15").Value = 0.01126098378240761
$ws.Range("C16").Value = 30.41673851013184
$ws.Range("D16").Value = -0.06326148986816804
$ws.Range("E16").Value = 0.004002016100340328
$ws.Range("C17").Value = 30.42951011657715
$ws.Range("D17").Value = -0.2604898834228493
$ws.Range("E17").Value = 0.06785497936564962
$ws.Range("C18").Value = 30.52583885192871
$ws.Range("D18").Value = -0.2241611480712891
$ws.Range("E18").Value = 0.05024822030463838
$ws.Range("C19").Value = 30.64208030700684
$ws.Range("D19").Value = -0.2979196929931618
$ws.Range("E19").Value = 0.08875614347313977
$ws.Range("C20").Value = 30.75008964538574
$ws.Range("D20").Value = -0.1999103546142607
$ws.Range("E20").Value = 0.03996414988199944
$ws.Range("C21").Value = 31.1014232635498
$ws.Range("D21").Value = 0.08142326354980867
$ws.Range("E21").Value = 0.006629747847101601
$ws.Range("C22").Value = 31.24580955505371
$ws.Range("D22").Value = 0.1258095550537064
$ws.Range("E22").Value = 0.01582804414281158
$ws.Range("C23").Value = 31.31645774841309
$ws.Range("D23").Value = 0.0364577484130848
$ws.Range("E23").Value = 0.001329167419351787
$ws.Range("C24").Value = 31.22794914245605
$ws.Range("D24").Value = -0.1520508575439408
$ws.Range("E24").Value = 0.02311946327984777
$ws.Range("C25").Value = 31.41512870788574
$ws.Range("D25").Value = -0.1648712921142561
$ws.Range("E25").Value = 0.02718254296342437
$ws.Range("B26").Value = 31.65000000000001
$ws.Range("C26").Value = 31.88678550720215
$ws.Range("D26").Value = 0.2367855072021428
$ws.Range("E26").Value = 0.056067376420976
$ws.Range("C27").Value = 32.44120407104492
$ws.Range("D27").Value = 0.5612040710449264
$ws.Range("E27").Value = 0.3149500093573988
$ws.Range("C28").Value = 32.40054702758789
$ws.Range("D28").Value = 0.1205470275878895
$ws.Range("E28").Value = 0.01453158586027539
$ws.Range("C29").Value = 32.50535202026367
$ws.Range("D29").Value = 0.05535202026366903
$ws.Range("E29").Value = 0.003063846147269627
$ws.Range("B30").Value = 32.84999999999999
$ws.Range("C30").Value = 32.76932907104492
$ws.Range("D30").Value = -0.08067092895507244
$ws.Range("E30").Value = 0.006507798778474345
$ws.Range("B31").Value = 32.90000000000001
$ws.Range("C31").Value = 32.95795059204102
$ws.Range("D31").Value = 0.05795059204100994
$ws.Range("E31").Value = 0.003358271117903565
$ws.Range("B32").Value = 33.09999999999999
$ws.Range("C32").Value = 32.89803695678711
$ws.Range("D32").Value = -0.2019630432128849
$ws.Range("E32").Value = 0.04078907082380963
$ws.Range("B33").Value = 33.40000000000001
$ws.Range("C33").Value = 33.65871047973633
$ws.Range("D33").Value = 0.2587104797363224
$ws.Range("E33").Value = 0.0669311123253981
$ws.Range("C34").Value = 33.65222549438477
$ws.Range("D34").Value = -0.04777450561523722
$ws.Range("E34").Value = 0.002282403386780332
$ws.Range("B35").Value = 34.09999999999999
$ws.Range("C35").Value = 33.83626174926758
$ws.Range("D35").Value = -0.2637382507324162
$ws.Range("E35").Value = 0.06955786489939483
$ws.Range("B36").Value = 34.40000000000001
$ws.Range("C36").Value = 34.3889045715332
$ws.Range("D36").Value = -0.01109542846680256
$ws.Range("E36").Value = 0.0001231085328619326
$ws.Range("B37").Value = 34.90000000000001
$ws.Range("C37").Value = 35.00925064086914
$ws.Range("D37").Value = 0.1092506408691349
$ws.Range("E37").Value = 0.0119357025303167
$ws.Range("C38").Value = 35.71824264526367
$ws.Range("D38").Value = 0.4182426452636747
$ws.Range("E38").Value = 0.174926910317156
$ws.Range("C39").Value = 35.99924087524414
$ws.Range("D39").Value = 0.2992408752441378
$ws.Range("E39").Value = 0.08954510141687763
$ws.Range("C40").Value = 36.01230239868164
$ws.Range("D40").Value = -0.2876976013183565
$ws.Range("E40").Value = 0.08276990980433602
$ws.Range("C41").Value = 36.59271621704102
$ws.Range("D41").Value = -0.2072837829589815
$ws.Range("E41").Value = 0.04296656667778617
$ws.Range("C42").Value = 37.18803024291992
$ws.Range("D42").Value = -0.1119697570800753
$ws.Range("E42").Value = 0.01253722650057107
$ws.Range("B43").Value = 37.90000000000001
$ws.Range("C43").Value = 37.9639778137207
$ws.Range("D43").Value = 0.06397781372069744
$ws.Range("E43").Value = 0.004093160648480262
$ws.Range("C44").Value = 38.41967391967773
$ws.Range("D44").Value = -0.08032608032226562
$ws.Range("E44").Value = 0.006452279179939069
$ws.Range("B45").Value = 38.90000000000001
$ws.Range("C45").Value = 39.0008659362793
$ws.Range("D45").Value = 0.1008659362792912
$ws.Range("E45").Value = 0.01017393710149803
$ws.Range("B46").Value = 39.40000000000001
$ws.Range("C46").Value = 39.5327262878418
$ws.Range("D46").Value = 0.1327262878417912
$ws.Range("E46").Value = 0.01761626748426201
$ws.Range("B47").Value = 39.90000000000001
$ws.Range("C47").Value = 39.54935836791992
$ws.Range("D47").Value = -0.3506416320800838
$ws.Range("E47").Value = 0.1229495541477849
$ws.Range("B48").Value = 40.09999999999999
$ws.Range("C48").Value = 40.00116348266602
$ws.Range("D48").Value = -0.09883651733397869
$ws.Range("E48").Value = 0.00976865715870987
$ws.Range("B49").Value = 40.59999999999999
$ws.Range("C49").Value = 40.52347183227539
$ws.Range("D49").Value = -0.07652816772460369
$ws.Range("E49").Value = 0.005856560455285074
$ws.Range("B50").Value = 40.90000000000001
$ws.Range("C50").Value = 40.82853317260742
$ws.Range("D50").Value = -0.07146682739258381
$ws.Range("E50").Value = 0.005107507417561367
$ws.Range("B51").Value = 41.20000000000001
$ws.Range("C51").Value = 41.49383163452148
$ws.Range("D51").Value = 0.2938316345214744
$ws.Range("E51").Value = 0.08633702944556132
$ws.Range("C52").Value = 0.1741580963134552
$ws.Range("E52").Value = 2.03359539015778
$ws.Range("E53").Value = 0.0406719078031556
